$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.302.53"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "1.840.01"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'239.12"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").Value = "'0.6247"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.07368"
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("D9").Value = "'0.2887"
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("D10").Value = "'24.71"
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("D11").Value = "'0.07717"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "1.828.77"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").Value = "'0.00001058"
$ws.Range("E14").Value = "  +3.35%  "
$ws.Range("D15").Value = "'0.6620"
$ws.Range("E15").Value = "  -2.97%  "
$ws.Range("D16").Value = "'81.37"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "'6.230"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").Value = "29.318.14"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "'236.33"
$ws.Range("E19").Value = "  +2.92%  "
$ws.Range("D20").Value = "'12.21"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'7.231"
$ws.Range("E22").Value = "  -3.83%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "'157.24"
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("D25").Value = "'8.406"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("D26").Value = "'0.1334"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("D28").Value = "'0.07067"
$ws.Range("E28").Value = "  +6.90%  "
$ws.Range("D29").Value = "'1.461"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'4.022"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.014"
$ws.Range("E32").Value = "  -1.98%  "
$ws.Range("D33").Value = "'1.152"
$ws.Range("D34").Value = "'1.777"
$ws.Range("E34").Value = "  -3.82%  "
$ws.Range("D35").Value = "'0.6840"
$ws.Range("E35").Value = "  -1.75%  "
$ws.Range("E36").Value = "  +0.54%  "
$ws.Range("D37").Value = "'0.01821"
$ws.Range("E37").Value = "  -2.60%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.779"
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.231.82"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("D40").Value = "'6.708"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").Value = "'0.9423"
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "1.999.02"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "'65.06"
$ws.Range("E46").Value = "  +7.79%  "
$ws.Range("D47").Value = "'6.933"
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("D48").Value = "'1.682"
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("D49").Value = "'8.860"
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("D50").Value = "'0.1127"
$ws.Range("D51").Value = "'0.3870"
$ws.Range("E51").Value = "  -1.52%  "
